$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 3 ("atopie (1)") to make room for
# the new "sexe (m)" row, pushing existing rows down.
$ws.Rows.Item(3).Insert()

# Fill the newly inserted row 3 with the "sexe (m)" data.
# Force text format on numeric-looking cells so they are stored as
# shared strings rather than being parsed as numbers.
$ws.Range("B3:D3").NumberFormat = "@"

$ws.Range("A3").Value = "sexe (m)"
$ws.Range("B3").Value = "3.281"
$ws.Range("C3").Value = "[0.666;24.445]"
$ws.Range("D3").Value = "0.176"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "-"
